$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60 (existing rows 60-86 shift down to 61-87)
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with the new record
$ws.Cells.Item(60, 1).Value = 8
$ws.Cells.Item(60, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(60, 3).Value = "Coquimbo"
$ws.Cells.Item(60, 4).Value = 44609
$ws.Cells.Item(60, 5).Value = 4
$ws.Cells.Item(60, 6).Value = "Fruta"
$ws.Cells.Item(60, 7).Value = 100109
$ws.Cells.Item(60, 8).Value = "Uva"
$ws.Cells.Item(60, 9).Value = 100109001
$ws.Cells.Item(60, 10).Value = "Uva"
$ws.Cells.Item(60, 11).Value = "Red Globe"
$ws.Cells.Item(60, 12).Value = "Primera"
$ws.Cells.Item(60, 13).Value = 500
$ws.Cells.Item(60, 14).Value = 11500
$ws.Cells.Item(60, 15).Value = 12000
$ws.Cells.Item(60, 16).Value = 11750
$ws.Cells.Item(60, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(60, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(60, 19).Value = 653
$ws.Cells.Item(60, 20).Value = 18
